$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.13
$ws.Range("C6").Value = -11.626
$ws.Range("D10").Value = -7.476000000000001
$ws.Range("A14").Value = -22.052
$ws.Range("B15").Value = 5.417999999999999
$ws.Range("D15").Value = -8.084999999999999
$ws.Range("A16").Value = -22.341
$ws.Range("C18").Value = -12.995
$ws.Range("D18").Value = -7.997
$ws.Range("C19").Value = -12.038
$ws.Range("A21").Value = -21.042
$ws.Range("B21").Value = 6.976000000000001
$ws.Range("D21").Value = -8.243
$ws.Range("B22").Value = 8.210000000000001
$ws.Range("D22").Value = -8.074999999999999
$ws.Range("A23").Value = -21.054
$ws.Range("B24").Value = 5.483000000000001
$ws.Range("D24").Value = -7.744999999999999
$ws.Range("A25").Value = -21.706
$ws.Range("A26").Value = -21.123
$ws.Range("B27").Value = 6.155000000000001
$ws.Range("B28").Value = 6.111000000000001
$ws.Range("A29").Value = -21.498
$ws.Range("D33").Value = -7.308
$ws.Range("C35").Value = -12.385
$ws.Range("B36").Value = 7.58
$ws.Range("B39").Value = 7.487
$ws.Range("A40").Value = -21.144
$ws.Range("C44").Value = -12.221
$ws.Range("B45").Value = 5.451000000000001
$ws.Range("D46").Value = -8.113
$ws.Range("C47").Value = -12.152
$ws.Range("B48").Value = 5.48
$ws.Range("B49").Value = 6.291000000000001
$ws.Range("D49").Value = -8.111999999999998
$ws.Range("C50").Value = -13.422
$ws.Range("C51").Value = -11.709
$ws.Range("B52").Value = 5.356
$ws.Range("C52").Value = -11.218
$ws.Range("A53").Value = -21.346
$ws.Range("B53").Value = 6.206999999999999
$ws.Range("B54").Value = 5.711
$ws.Range("C55").Value = -13.445
$ws.Range("D56").Value = -7.673
$ws.Range("A57").Value = -21.682
$ws.Range("B57").Value = 6.129
$ws.Range("C57").Value = -12.988
$ws.Range("C58").Value = -12.943
$ws.Range("A59").Value = -22.396
$ws.Range("D61").Value = -8.032
$ws.Range("C64").Value = -10.967
$ws.Range("A65").Value = -21.728
$ws.Range("D65").Value = -7.328
$ws.Range("C66").Value = -11.675
$ws.Range("D66").Value = -7.556
$ws.Range("A69").Value = -21.555
$ws.Range("B70").Value = 4.760000000000001
$ws.Range("B71").Value = 4.972
$ws.Range("D74").Value = -8.081
$ws.Range("D75").Value = -7.621
$ws.Range("D77").Value = -7.744
$ws.Range("A79").Value = -21.084
$ws.Range("C80").Value = -12.108
$ws.Range("A83").Value = -22.022
$ws.Range("C83").Value = -12.615
$ws.Range("B86").Value = 5.555000000000001
$ws.Range("B87").Value = 4.834999999999999
$ws.Range("D87").Value = -8.378
$ws.Range("D88").Value = -7.605000000000001
$ws.Range("B89").Value = 5.747999999999999
$ws.Range("A91").Value = -21.173
$ws.Range("C92").Value = -10.97
$ws.Range("A93").Value = -21.628
$ws.Range("C94").Value = -11.935
$ws.Range("C96").Value = -11.655
$ws.Range("C97").Value = -11.155
$ws.Range("A100").Value = -22.043
$ws.Range("D100").Value = -7.722
$ws.Range("B101").Value = 5.436
$ws.Range("C101").Value = -11.619
$ws.Range("D101").Value = -7.6
$ws.Range("A103").Value = -21.969
